$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# C2 used to hold "Login and Logout" -> now holds "Twitter Login"
$ws.Range("C2").Value = "Twitter Login"

# D2 used to hold "No" -> now holds "Yes"
$ws.Range("D2").Value = "Yes"

# D4 used to hold "Yes" -> now holds "No"
$ws.Range("D4").Value = "No"

# Update the selected cell to match the new active selection (C2)
$ws.Range("C2").Select()
